$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

$ws.Range("B2").Value = "LG"
$ws.Range("B3").Value = "LG.5M"
